$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.625.47"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "2.461.40"
$ws.Range("E3").Value = "  -0.97%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.04"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.36"
$ws.Range("E6").Value = "  -0.82%  "

# Row 7
$ws.Range("E7").Value = "  +0.61%  "

# Row 8
$ws.Range("E8").Value = "  +0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  +3.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("E10").Value = "  +0.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  +2.20%  "

# Row 12
$ws.Range("E12").Value = "  +0.63%  "

# Row 13
$ws.Range("D13").Value = "2.838.51"
$ws.Range("E13").Value = "  -0.90%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  +0.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.81"
$ws.Range("E15").Value = "  +3.37%  "

# Row 16
$ws.Range("D16").Value = "2.464.28"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").Value = "  +2.53%  "

# Row 18
$ws.Range("D18").Value = "41.610.43"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  +2.82%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  +1.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.66"
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.31"
$ws.Range("E22").Value = "  +2.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.21"
$ws.Range("E23").Value = "  +1.12%  "

# Row 24
$ws.Range("E24").Value = "  +0.54%  "

# Row 25
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.29"
$ws.Range("E27").Value = "  -0.63%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  +0.63%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.71"
$ws.Range("E29").Value = "  +1.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.07"
$ws.Range("E30").Value = "  -2.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.04"
$ws.Range("E31").Value = "  +1.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.46"
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("E33").Value = "  +0.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0760"
$ws.Range("E34").Value = "  +0.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  -1.31%  "

# Row 36
$ws.Range("E36").Value = "  -2.77%  "

# Row 37
$ws.Range("E37").Value = "  -2.52%  "

# Row 38
$ws.Range("E38").Value = "  +1.70%  "

# Row 39
$ws.Range("E39").Value = "  +2.07%  "

# Row 40
$ws.Range("E40").Value = "  -1.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.94"
$ws.Range("E41").Value = "  -3.92%  "

# Row 42
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("D43").Value = "1.972.61"
$ws.Range("E43").Value = "  +1.53%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.91"
$ws.Range("E44").Value = "  -2.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  -1.27%  "

# Row 47
$ws.Range("E47").Value = "  +2.72%  "

# Row 48
$ws.Range("D48").Value = "2.695.29"
$ws.Range("E48").Value = "  -1.15%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.83"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.81"
$ws.Range("E50").Value = "  -0.15%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.24"
$ws.Range("E51").Value = "  +4.07%  "
